# Update cleaned SR time series with new Hucuktlis CVs

$wb = $excel.ActiveWorkbook

# --- 1. Update metadata definition text for H_cv (sheet "metadata", cell B17) ---
$wsMeta = $wb.Worksheets.Item("metadata")
$wsMeta.Range("B17").Value = "Coefficient of variation on harvest data. Historical (prior to 2011) Hucuktlis Sockeye harvest rate predictions were derived from a linear model. Harvest data for Somass and Hucuktlis post-2011 are assumed to be precise."

# --- 2. Update S_cv values (sheet "S-R data", column Q) for HUC rows ---
$wsData = $wb.Worksheets.Item("S-R data")

$updates = @{
    59  = 0.2496267659213876
    60  = 0.2111714885758418
    61  = 0.2922206205809506
    62  = 0.2565178192417154
    63  = 0.2451098642822979
    64  = 0.2419729833867256
    65  = 0.2765457477820611
    66  = 0.231004202683785
    67  = 0.2464772207893425
    68  = 0.2480305876316884
    69  = 0.2553852420154066
    70  = 0.2550631891375683
    71  = 0.2553581160154582
    72  = 0.2553220566127163
    73  = 0.2821999907530116
    74  = 0.2525725048464748
    75  = 0.178692682150781
    76  = 0.2562130945777705
    77  = 0.2620673654633024
    78  = 0.2450461298290044
    85  = 0.2487218717321054
    86  = 0.1827938490156061
    89  = 0.1276418702964884
    90  = 0.2288747452863037
    91  = 0.2903052391163663
    92  = 0.2556215086392382
    93  = 0.241441424039182
    94  = 0.2497769802817572
    95  = 0.2440455592874953
    96  = 0.2425905410682488
    97  = 0.2388851566095718
    98  = 0.2506997776881581
    99  = 0.1397559218405289
    100 = 0.1759984840093463
    101 = 0.1761708452318516
    102 = 0.1677283430730512
}

foreach ($row in $updates.Keys) {
    $wsData.Range("Q$row").Value = $updates[$row]
}
